$d = $word.ActiveDocument

# Locate the final paragraph of the document ("All by all Im very happy about my app.")
$lastParaIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($lastParaIndex)
$targetRange = $targetPara.Range

# Build the replacement OOXML package fragment: the original sentence (bookmark
# removed here; it will be re-added further down, attached to the new final
# "At last I added..." paragraph) followed by all of the new diary paragraphs,
# ending with a trailing blank paragraph.
$bodyInner = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>All by all Im very happy about my app.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(26 Juni 2018)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Today I m going to fix all my bugs. One bug that was really giving me headache was the fact that somehow in my ActivitiesTableViewController my multidimensional array kept getting bigger when I liked an activity and and reloaded my TableView. And this resulted in a section error. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I used several print functions to fix this problem and I found out that the fetchActivitiyId function triggerd more often than the fetchActivities function which resulted in that the fetchPartActivities function was also called more often than fetchActivities, because I call fetchPartActivities everytime at the end of fetchActivityId. But because is only set my allActivities multi d array to [] in fetchActivities, my array kept growing everytime. So my solution was to directly after I append activities in FetchActivities to also append a empty array [], and set this allActivities[1] = [] everytime when I was going to set new values for partActivities. So in this way it wasn’t becoming bigger but just rewriting the information.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Also I fixed the friendsDetailviewcontroller from bugs. It was first not showing the friendsParticipating activities and when I thought I fixed it.. it was crashing because finding nil when unwrapping optional value. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">This is because I forgot to call ref= Databas… in viewDidLoad() </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">I have had this problem earlier before. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>So I almost fixed all my bugs and I only need to fix the delete function IF I am going to do that…</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>But I learned how to create my own buttons today with Sketch. I have seen several master tutorials a</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>nd I am proud that I succeeded.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">As well I changed the profile images to a complete round format, which already makes my design much user friendly in my opinion. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I did try to change my textfields as well, with a custom textfield from online that I downloaded via terminal with pod install, but my app doesn’t want to autocomplete after I use that class of the custom textfield. So eventually I just customized the textfields myself.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>And also made an Icon Logo for the app, I had to make several formats, but its done.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>At last I added some shadows to my swipe function which makes this even more sleek.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replacing the whole paragraph range (text + end-of-paragraph mark) inserts the
# new content immediately before it and leaves the original paragraph in place,
# emptied of its runs (and therefore of the old _GoBack bookmark). That emptied
# paragraph becomes the new last paragraph of the document, ready to receive the
# final "Tomorrow I want..." sentence below.
[void]$targetRange.InsertXML($packageXml)

# Fill in the text of the now-empty trailing paragraph.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.InsertAfter('Tomorrow I want to have look at the proper design format for apps. And try to improve the quality even more.')

# The freshly typed run does not automatically inherit the paragraph mark's
# run formatting, so restore the en-US language explicitly to match the rest
# of the document.
$finalRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$finalRange.LanguageID = "en-US"
